$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.406.35'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -3.48%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.669.11'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.29%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '218.98'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.21%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.5173'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.88%  '
$ws.Range("E7").Value = '  +0.37%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.06461'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.89%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2573'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("E10").Value = '  -3.99%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07668'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.683.83'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.352'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.73%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '1.897.70'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.5549'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.02%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0₅8058'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.38%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '64.79'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -4.33%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '26.437.62'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -3.34%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.44%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '210.75'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.46%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.426'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("E22").Value = '  -2.96%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.901'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.22%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.34%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '144.69'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.94%  '
$ws.Range("E26").Value = '  -1.98%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.1169'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.75%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.010'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.57%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '15.83'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.97%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05264'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.85%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.263'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.382'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.48%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.228'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -5.80%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.577'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.14%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.764'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.88%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.377'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.22%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9305'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.93%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.5734'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.19%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.155.47'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +10.66%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.01606'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.49%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.8531'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  +0.40%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '5.659'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.40%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '100.26'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.807.10'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0₈111'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.4496'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.21%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '56.09'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.32%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.946'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.64%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05112'
$c.Style = "Normal"
